# Rotate the comma-separated "Recorded By" list in column G one position to
# the right (the last entry in the list is moved to the front) for every row
# that has more than one entry. Single-entry cells are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $raw = $cell.Value2

    if ($raw -ne $null -and $raw.ToString().Contains(",")) {
        $parts = $raw.ToString().Split(",")

        $trimmed = @()
        foreach ($p in $parts) {
            $trimmed += $p.Trim()
        }

        $n = $trimmed.Length
        $last = $trimmed[$n - 1]
        $rest = $trimmed[0..($n - 2)]
        $rotated = @($last) + $rest

        $cell.Value = $rotated -join ", "
    }
}
